# Reorder the "group-code / group-name / category-name" columns (E, F, G)
# on the SectorGroup sheet so that the layout becomes:
#   E = codeforiati:category-name
#   F = codeforiati:group-code
#   G = codeforiati:group-name
# (previously: E = group-code, F = group-name, G = category-name)
#
# For every data row this is a simple right-rotation of the three cell
# values: newE = oldG, newF = oldE, newG = oldF.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 1 (E1:G1)
$ws.Range("E1").Value = "codeforiati:category-name"
$ws.Range("F1").Value = "codeforiati:group-code"
$ws.Range("G1").Value = "codeforiati:group-name"

# Determine last used row based on column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $oldE = $ws.Cells.Item($r, 5).Value2
    $oldF = $ws.Cells.Item($r, 6).Value2
    $oldG = $ws.Cells.Item($r, 7).Value2

    $ws.Cells.Item($r, 5).Value = $oldG
    $ws.Cells.Item($r, 6).Value = $oldE
    $ws.Cells.Item($r, 7).Value = $oldF
}
